# Reorder the 5 news-article records (rows 2-6) of the historical-distance
# sheet. Each record's title / timestamp / uri travel together; only the
# row position of each record changes. New top-to-bottom order:
#   1. Live results from the North Carolina primary
#   2. NBC's Final Battleground Map Shows Clinton With a Significant Lead
#   3. Clinton's North Carolina Firewall
#   4. The Final 15: The Latest Polls in the Swing States That Will Decide
#      the Election
#   5. NC approves 27 candidates for presidential primary ballots

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
    "Live results from the North Carolina primary",
    "NBC's Final Battleground Map Shows Clinton With a Significant Lead",
    "Clinton’s North Carolina Firewall",
    "The Final 15: The Latest Polls in the Swing States That Will Decide the Election",
    "NC approves 27 candidates for presidential primary ballots"
)

$timestamps = @(
    "2016-03-15T00:00:00UTC",
    "2016-11-07T13:17:36UTC",
    "2016-11-04T15:03:00UTC",
    "2016-11-07T15:39:00UTC",
    "2015-04-13T00:00:00UTC"
)

$uris = @(
    "http://graphics.latimes.com/election-2016-north-carolina-results/",
    "http://www.nbcnews.com/storyline/2016-election-day/nbc-s-final-battleground-map-shows-clinton-edge-n678926",
    "https://www.usnews.com/news/the-run-2016/articles/2016-11-04/hillary-clintons-north-carolina-firewall-vs-donald-trump",
    "https://abcnews.go.com/Politics/final-15-latest-polls-swing-states-decide-election/story?id=43277505",
    "http://www.wral.com/board-of-elections-sets-presidential-primary-ballots/15215886/"
)

# Write the reordered title/timestamp values.
for ($i = 0; $i -lt $titles.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $titles[$i]
    $ws.Cells.Item($r, 2).Value = $timestamps[$i]
}

# Hyperlinks travel with their uri, so rebuild the whole collection in the
# new row order (this engine's per-cell Hyperlinks.Delete() clears the
# entire sheet, so do the delete once up front).
$ws.Hyperlinks.Delete()
for ($i = 0; $i -lt $uris.Length; $i++) {
    $r = $i + 2
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $uris[$i]
    $ws.Hyperlinks.Add($cell, $uris[$i]) | Out-Null
    $cell.Style = "Hyperlink"
}
